$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "This is 1st change"

$ws.Range("A2").Select()
